$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates (odds refresh)
$ws.Range("M2").Value = 1.08
$ws.Range("N2").Value = 8

# Row 3 updates (odds refresh)
$ws.Range("G3").Value = 1.7
$ws.Range("M3").Value = 1.08
$ws.Range("N3").Value = 8
$ws.Range("AF3").Value = 81
$ws.Range("AO3").Value = 9
$ws.Range("AQ3").Value = 29
$ws.Range("AZ3").Value = 126

# Row 4 updates (odds refresh)
$ws.Range("G4").Value = 4
$ws.Range("I4").Value = 2.05
$ws.Range("J4").Value = 4.5
$ws.Range("M4").Value = 1.11
$ws.Range("N4").Value = 6.5
$ws.Range("Q4").Value = 2.5
$ws.Range("R4").Value = 1.5
$ws.Range("U4").Value = 2.1
$ws.Range("V4").Value = 1.67
$ws.Range("W4").Value = 9
$ws.Range("AC4").Value = 6.5
$ws.Range("AD4").Value = 6
$ws.Range("AH4").Value = 6
$ws.Range("AK4").Value = 19
$ws.Range("AW4").Value = 4

# Row 6 updates (odds refresh)
$ws.Range("N6").Value = 7.5
$ws.Range("W6").Value = 12
$ws.Range("AA6").Value = 51
$ws.Range("AD6").Value = 7
$ws.Range("AK6").Value = 12
$ws.Range("AN6").Value = 7
$ws.Range("AO6").Value = 34

# Row 7 updates (odds refresh)
$ws.Range("G7").Value = 1.75
$ws.Range("M7").Value = 1.1
$ws.Range("N7").Value = 7
$ws.Range("W7").Value = 5
$ws.Range("AC7").Value = 7
$ws.Range("AN7").Value = 3.5

# Row 8 updates (odds refresh)
$ws.Range("G8").Value = 2.6
$ws.Range("I8").Value = 2.37
$ws.Range("J8").Value = 3.1
$ws.Range("L8").Value = 2.85
$ws.Range("W8").Value = 11.5
$ws.Range("X8").Value = 15.5
$ws.Range("Z8").Value = 30
$ws.Range("AE8").Value = 11.5
$ws.Range("AH8").Value = 11.25
$ws.Range("AJ8").Value = 9.25
$ws.Range("AN8").Value = 4.8
$ws.Range("AO8").Value = 13.5
$ws.Range("AQ8").Value = 55
$ws.Range("AX8").Value = 11.75
$ws.Range("AY8").Value = 16.5
$ws.Range("BA8").Value = 60

# Row 13 updates (odds refresh)
$ws.Range("O13").Value = 1.29
$ws.Range("P13").Value = 3.5
$ws.Range("Q13").Value = 1.95
$ws.Range("R13").Value = 1.9

# Row 14 updates (odds refresh)
$ws.Range("G14").Value = 2.45
$ws.Range("I14").Value = 2.67
$ws.Range("J14").Value = 3.1
$ws.Range("K14").Value = 2.05
$ws.Range("L14").Value = 3.3
$ws.Range("S14").Value = 1.44
$ws.Range("T14").Value = 2.62
$ws.Range("V14").Value = 1.88
$ws.Range("W14").Value = 7.6
$ws.Range("X14").Value = 11.75
$ws.Range("Z14").Value = 27
$ws.Range("AA14").Value = 22
$ws.Range("AD14").Value = 6.1
$ws.Range("AG14").Value = 600
$ws.Range("AH14").Value = 7.9
$ws.Range("AI14").Value = 13
$ws.Range("AJ14").Value = 10.25
$ws.Range("AK14").Value = 30
$ws.Range("AL14").Value = 24
$ws.Range("AM14").Value = 35
$ws.Range("AN14").Value = 4.35
$ws.Range("AO14").Value = 13.5
$ws.Range("AP14").Value = 22
$ws.Range("AQ14").Value = 60
$ws.Range("AR14").Value = 100
$ws.Range("AT14").Value = 2.62
$ws.Range("AW14").Value = 4.6
$ws.Range("AX14").Value = 15
$ws.Range("AZ14").Value = 70

